# Apply crypto price/volume updates from the Oct 25 2024 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.039.43"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.18%  "

# Row 3: update D3, E3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.543.47"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.61%  "

# Row 4: update D4, E4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5: update D5, E5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.86"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.34%  "

# Row 6: update D6, E6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.57"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.21%  "

# Row 7: update E7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.09%  "

# Row 8: update D8, E8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.527"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.21%  "

# Row 9: update D9, E9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.541.95"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.56%  "

# Row 10: update E10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.49%  "

# Row 11: update E11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.86%  "

# Row 12: update D12, E12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.06"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.83%  "

# Row 13: update E13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.02%  "

# Row 14: update D14, E14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.53"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.23%  "

# Row 15: update D15, E15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.002.27"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.44%  "

# Row 16: update E16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.99%  "

# Row 17: update D17, E17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.898.72"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.30%  "

# Row 18: update D18, E18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.37"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +137.61%  "

# Row 19: update D19, E19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.527.95"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.01%  "

# Row 20: update D20, E20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.78"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.39%  "

# Row 21: update D21, E21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.99"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.96%  "

# Row 22: update D22, E22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "370.02"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.27%  "

# Row 23: update D23, E23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.16"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.58%  "

# Row 24: update D24, E24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.59"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.34%  "

# Row 25: update D25, E25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.64"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.64%  "

# Row 26: update B26, C26, D26, E26
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "SuiNetwork"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.94"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.12%  "

# Row 27: update B27, C27, D27, E27
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.06%  "

# Row 28: update D28, E28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.98"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.56%  "

# Row 29: update D29, E29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.609.22"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.76%  "

# Row 30: update D30, E30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0974"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.67%  "

# Row 31: update D31, E31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.49"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.76%  "

# Row 32: update D32, E32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "542.01"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.36%  "

# Row 33: update E33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.72%  "

# Row 34: update E34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.28%  "

# Row 35: update D35, E35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.130"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.53%  "

# Row 36: update D36, E36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.02%  "

# Row 37: update D37, E37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.48"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.58%  "

# Row 38: update E38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.65%  "

# Row 39: update D39, E39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.18"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.83%  "

# Row 40: update E40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.95%  "

# Row 41: update D41, E41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.18"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.73%  "

# Row 42: update E42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.56%  "

# Row 43: update D43, E43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.79"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.64%  "

# Row 44: update D44, E44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.59"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.10%  "

# Row 45: update E45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.06%  "

# Row 46: update D46, E46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.30"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.96%  "

# Row 47: update D47, E47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0290"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.69%  "

# Row 48: update D48, E48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "147.88"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.81%  "

# Row 49: update E49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.77%  "

# Row 50: update D50, E50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.554"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.60%  "

# Row 51: update D51, E51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.72"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.53%  "
